$d = $word.ActiveDocument

# This document has 4 M2Doc field codes (complex Word fields whose
# instruction text holds the M2Doc template syntax, e.g. " m:v.name ").
# This edit rewrites each of them from the Word-field representation
# (fldChar begin / instrText* / fldChar end) into plain text runs using
# the M2Doc "{m:...}" textual token syntax, stripping the single
# leading/trailing space that Word requires around field instructions.
#
# Fields are always addressed through the document-level `$d.Fields`
# collection (field 1, re-fetched each time) since deleting a field
# shifts every later field's index down by one, and `Range.Fields` on
# a single paragraph's Range is not reliable in this runtime.
#
# For each field, the insertion point for the replacement text is
# recorded *before* deletion as "where the field's instruction text
# starts, minus one" (that extra character is the fldChar "begin" run
# immediately preceding the instruction text) -- i.e. exactly where the
# field as a whole begins. Because deleting the field only removes
# characters at/after that position, the position itself still points
# at the right spot afterwards, regardless of what text follows the
# field in the same paragraph (e.g. the trailing "," after the second
# field).

function Convert-FieldToText([string]$newText) {
    $f = $d.Fields.Item(1)
    $fieldStart = $f.Code.Start - 1
    $f.Delete()
    $rng = $d.Range($fieldStart, $fieldStart)
    $rng.InsertBefore($newText)
}

# Field 1 : " for v | self.eClassifiers" -> "{m:for v | self.eClassifiers}"
Convert-FieldToText "{m:for v | self.eClassifiers}"

# Field 2 : " m:v.name " -> "{m:v.name}" (the trailing "," stays put)
Convert-FieldToText "{m:v.name}"

# Field 3 : " m:'newParagraph'.asPagination() " -> "{m:'newParagraph'.asPagination()}"
# This field also wraps a _GoBack bookmark (right after "newParagraph")
# which must be preserved in the rewritten text.
$f3 = $d.Fields.Item(1)
$fieldStart3 = $f3.Code.Start - 1
$f3.Delete()
$rng3 = $d.Range($fieldStart3, $fieldStart3)
$rng3.InsertBefore("{m:'newParagraph'.asPagination()}")
$bmPos = $fieldStart3 + "{m:'newParagraph".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Field 4 : " m:endfor " -> "{m:endfor}"
Convert-FieldToText "{m:endfor}"
